$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update employee records (rows 2-4). Columns: A firstName, B middleName,
# C lastName, D photo, E username, F password, G confirmPassword.
$ws.Range("A2").Value = "Bertokr"
$ws.Range("C2").Value = "Sapirtom"
$ws.Range("E2").Value = "Bertokk325"

$ws.Range("A3").Value = "Weportt"
$ws.Range("C3").Value = "Derakoll"
$ws.Range("E3").Value = "Weportt325"

$ws.Range("A4").Value = "Xeelopp"
$ws.Range("C4").Value = "Pomedorr"
$ws.Range("E4").Value = "Xeelopp325"

# Update the active cell selection shown when the workbook is opened.
$ws.Range("E7").Select()
